$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell BQ1 (date label), matching the style/format of BP1.
$ws.Range("BQ1").Value = "'2023-04-28"
$ws.Range("BP1").Copy()
$ws.Range("BQ1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update existing BP column values and add new BQ column values for rows 2-23.
$ws.Range("BP2").Value = 802297
$ws.Range("BQ2").Value = 817325

$ws.Range("BP3").Value = 994148
$ws.Range("BQ3").Value = 1003178

$ws.Range("BP4").Value = 1939992
$ws.Range("BQ4").Value = 2282896

$ws.Range("BP5").Value = 1452796
$ws.Range("BQ5").Value = 1426117

$ws.Range("BP6").Value = 4656640
$ws.Range("BQ6").Value = 4797634

$ws.Range("BP7").Value = 1692118
$ws.Range("BQ7").Value = 1696101

$ws.Range("BP8").Value = 11355904
$ws.Range("BQ8").Value = 11466746

$ws.Range("BP9").Value = 2714332
$ws.Range("BQ9").Value = 2945892

$ws.Range("BP10").Value = 1155275
$ws.Range("BQ10").Value = 1227897

$ws.Range("BP11").Value = 1674776
$ws.Range("BQ11").Value = 1724755

$ws.Range("BP12").Value = 2215071
$ws.Range("BQ12").Value = 2359622

$ws.Range("BP13").Value = 1146497
$ws.Range("BQ13").Value = 1225727

$ws.Range("BP14").Value = 1250117
$ws.Range("BQ14").Value = 1274353

$ws.Range("BP15").Value = 9964544
$ws.Range("BQ15").Value = 10020517

$ws.Range("BP16").Value = 1285398
$ws.Range("BQ16").Value = 1295615

$ws.Range("BP17").Value = 528114
$ws.Range("BQ17").Value = 539895

$ws.Range("BP18").Value = 1740493
$ws.Range("BQ18").Value = 1792203

$ws.Range("BP19").Value = 603685
$ws.Range("BQ19").Value = 633193

$ws.Range("BP20").Value = 3471629
$ws.Range("BQ20").Value = 3497906

$ws.Range("BP21").Value = 2397168
$ws.Range("BQ21").Value = 2943599

$ws.Range("BP22").Value = 421279
$ws.Range("BQ22").Value = 425913

$ws.Range("BP23").Value = 10273671
$ws.Range("BQ23").Value = 10295867
